# Update NATMI ligand/receptor expression statistics with newly recomputed
# TPM-based values. The workbook lists, for every (Sending cluster, Target
# cluster) pair, the ligand stats for the sending cluster (columns G-J), the
# receptor stats for the target cluster (columns M-P), and edge weights that
# are simply the products of the ligand/receptor average & total values and
# their specificities (columns Q-T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand (sending-cluster) statistics:
#   Ligand average expr, Ligand total expr, Ligand avg specificity, Ligand total specificity
$ligandStats = @{
    "ECs"              = @(29.20950566666667, 87.628517,       0.01829497698069002, 0.01840828041918582)
    "FAPs"             = @(1458.280985666667, 4374.842957,     0.913374480506715,   0.9190311407684336)
    "Inflammatory-Mac" = @(57.98602933333333, 173.958088,      0.03631876156896331, 0.03654368891224535)
    "MuSCs"            = @(29.481085,         58.96217,        0.01846507700595112, 0.01238628926567028)
    "Resolving-Mac"    = @(21.628479,         64.885437,       0.01354670393768061, 0.01363060063446486)
}

# New receptor (target-cluster) statistics:
#   Receptor average expr, Receptor total expr, Receptor avg specificity, Receptor total specificity
$receptorStats = @{
    "ECs"              = @(16.14072933333334, 48.42218800000001, 0.03423048004954622, 0.03634868370049611)
    "FAPs"             = @(81.06331633333333, 243.189949,        0.1719151703242873,  0.1825533892714798)
    "Inflammatory-Mac" = @(168.70371,         506.11113,         0.3577786889414888,  0.3799182594076638)
    "MuSCs"            = @(82.43477250000001, 164.869545,        0.1748236883957081,  0.1237612588479007)
    "Resolving-Mac"    = @(123.1883796666667, 369.565139,        0.2612519722889696,  0.2774184087724594)
}

for ($r = 2; $r -le 26; $r++) {
    $sending = $ws.Cells.Item($r, 1).Text   # column A: Sending cluster
    $target  = $ws.Cells.Item($r, 4).Text   # column D: Target cluster

    $lig = $ligandStats[$sending]
    $rec = $receptorStats[$target]

    $gVal = $lig[0]
    $hVal = $lig[1]
    $iVal = $lig[2]
    $jVal = $lig[3]

    $mVal = $rec[0]
    $nVal = $rec[1]
    $oVal = $rec[2]
    $pVal = $rec[3]

    $ws.Cells.Item($r, 7).Value  = $gVal   # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $hVal   # H: Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $iVal   # I: Ligand derived specificity of average expression value
    $ws.Cells.Item($r, 10).Value = $jVal   # J: Ligand derived specificity of total expression value

    $ws.Cells.Item($r, 13).Value = $mVal   # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $nVal   # N: Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $oVal   # O: Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value = $pVal   # P: Receptor derived specificity of total expression value

    $ws.Cells.Item($r, 17).Value = $gVal * $mVal   # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $hVal * $nVal   # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $iVal * $oVal   # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $jVal * $pVal   # T: Edge total expression derived specificity
}
